$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The recorded edit simply re-keys A1 from 13 to 10.
$ws.Range("A1").Value = 10

# Author's resave also left the view parked back on A1 (the implicit
# default) instead of the previously selected C1.
$ws.Range("A1").Select() | Out-Null
